# IC Module: bug fixed, topic speaker list not loaded
$d = $word.ActiveDocument

# 1) Replace "Председатель Комитета" with "HEAD" (keep same run's formatting, but language will be set after)
$d.Content.Find.Execute("Председатель Комитета", $false, $false, $false, $false, $false,
                         $true, 1, $false, "HEAD", 2)

# 2) Replace "       Таджияков Г.Б." (with its leading tab run / spaces) with "NAME"
$d.Content.Find.Execute("Таджияков Г.Б.", $false, $false, $false, $false, $false,
                         $true, 1, $false, "NAME", 2)

Write-Host "Done text replace"
